$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $CellAddr, $TextVal)
    $rng = $Worksheet.Range($CellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $TextVal
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "26.421.89"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Value = "1.689.30"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  +0.52%  "
Set-TextValue $ws "D5" "219.02"
$ws.Range("E5").Value = "  +1.20%  "
Set-TextValue $ws "D6" "0.5530"
$ws.Range("E6").Value = "  +8.59%  "
Set-TextValue $ws "D8" "0.2713"
$ws.Range("E8").Value = "  +2.40%  "
Set-TextValue $ws "D9" "0.06496"
$ws.Range("E9").Value = "  +1.72%  "
Set-TextValue $ws "D10" "22.17"
$ws.Range("E10").Value = "  +1.77%  "
Set-TextValue $ws "D11" "0.07597"
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws "D12" "4.567"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.682.24"
$ws.Range("E13").Value = "  +0.56%  "
Set-TextValue $ws "D14" "0.5829"
$ws.Range("E14").Value = "  +0.16%  "
Set-TextValue $ws "D15" "0.000008490"
$ws.Range("E15").Value = "  -0.43%  "
Set-TextValue $ws "D16" "65.54"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").Value = "26.468.97"
$ws.Range("E17").Value = "  +1.37%  "
Set-TextValue $ws "D18" "4.952"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("E19").Value = "  +0.50%  "
Set-TextValue $ws "D20" "10.98"
$ws.Range("E20").Value = "  +1.96%  "
Set-TextValue $ws "D21" "191.15"
$ws.Range("E21").Value = "  -0.01%  "
Set-TextValue $ws "D22" "6.263"
Set-TextValue $ws "D23" "1.010"
$ws.Range("E23").Value = "  +0.50%  "
Set-TextValue $ws "D24" "149.92"
$ws.Range("E24").Value = "  +3.83%  "
Set-TextValue $ws "D25" "0.1325"
$ws.Range("E25").Value = "  +10.47%  "
Set-TextValue $ws "D26" "7.926"
$ws.Range("E26").Value = "  +4.35%  "
Set-TextValue $ws "D27" "15.88"
$ws.Range("E27").Value = "  +1.45%  "
Set-TextValue $ws "D28" "0.06349"
$ws.Range("E28").Value = "  -4.29%  "
Set-TextValue $ws "D29" "1.405"
$ws.Range("E29").Value = "  +5.24%  "
Set-TextValue $ws "D30" "1.328"
$ws.Range("E30").Value = "  +0.94%  "
Set-TextValue $ws "D31" "3.594"
$ws.Range("E31").Value = "  +1.34%  "
Set-TextValue $ws "D32" "3.589"
$ws.Range("E32").Value = "  +2.34%  "
Set-TextValue $ws "D33" "1.679"
$ws.Range("E33").Value = "  +1.40%  "
Set-TextValue $ws "D34" "1.045"
$ws.Range("E34").Value = "  +2.75%  "
Set-TextValue $ws "D35" "0.6256"
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("E36").Value = "  +1.63%  "
Set-TextValue $ws "D37" "2.723"
$ws.Range("E37").Value = "  +1.54%  "
Set-TextValue $ws "D38" "6.252"
$ws.Range("E38").Value = "  -0.95%  "
Set-TextValue $ws "D39" "0.01640"
$ws.Range("E39").Value = "  +2.85%  "
$ws.Range("D40").Value = "1.118.95"
$ws.Range("E40").Value = "  +2.21%  "
Set-TextValue $ws "D41" "0.8798"
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("E42").Value = "  +0.61%  "
Set-TextValue $ws "D43" "100.88"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "1.838.70"
$ws.Range("E44").Value = "  +1.32%  "
Set-TextValue $ws "D45" "0.00000000111"
$ws.Range("E45").Value = "  -0.67%  "
Set-TextValue $ws "D46" "57.56"
$ws.Range("E46").Value = "  +2.13%  "
Set-TextValue $ws "D47" "8.223"
$ws.Range("E47").Value = "  +2.15%  "
Set-TextValue $ws "D48" "1.007"
$ws.Range("E48").Value = "  +0.08%  "
Set-TextValue $ws "D49" "0.05286"
$ws.Range("E49").Value = "  +1.18%  "
Set-TextValue $ws "D50" "0.4303"
$ws.Range("E50").Value = "  +0.39%  "
Set-TextValue $ws "D51" "6.101"
$ws.Range("E51").Value = "  +1.24%  "
